# Natmi following Dr Hou advice
# Re-run of the NATMI LR-pair analysis (Thbs1-Itga3) now includes the
# previously-missing "M2" target cluster, expanding each sending cluster's
# block from 3 target rows to 4, and adds the new "sCs" sending-cluster
# block (rows 14-17). Existing rows are refreshed with the recomputed
# statistics from the new run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full refreshed data for rows 2-17 (columns A-T), in row order:
# Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# Ligand-expressing cells, Ligand detection rate, Ligand average expression value,
# Ligand total expression value, Ligand derived specificity of average expression value,
# Ligand derived specificity of total expression value, Receptor-expressing cells,
# Receptor detection rate, Receptor average expression value, Receptor total expression value,
# Receptor derived specificity of average expression value, Receptor derived specificity of total expression value,
# Edge average expression weight, Edge total expression weight,
# Edge average expression derived specificity, Edge total expression derived specificity
$data = @(
    @("ECs","Thbs1","Itga3","ECs",3,1,18.82499266666667,56.474978,0.06886869772378311,0.0688686977237831,3,1,12.77892533333333,38.336776,0.7206984557633653,0.7206984557633654,240.5631756878809,2165.068581190928,0.04963356409996448,0.04963356409996448),
    @("ECs","Thbs1","Itga3","FAPs",3,1,18.82499266666667,56.474978,0.06886869772378311,0.0688686977237831,3,1,0.9898276666666668,2.969483,0.05582372947885773,0.05582372947885774,18.63349856626378,167.701487096374,0.003844507551293693,0.003844507551293693),
    @("ECs","Thbs1","Itga3","M2",3,1,18.82499266666667,56.474978,0.06886869772378311,0.0688686977237831,1,0.3333333333333333,0.02976266666666667,0.089288,0.001678537697541373,0.001678537697541373,0.5602819817404445,5.042537835664,0.0001155987053099517,0.0001155987053099517),
    @("ECs","Thbs1","Itga3","sCs",3,1,18.82499266666667,56.474978,0.06886869772378311,0.0688686977237831,3,1,3.932791,11.798373,0.2217992770602354,0.2217992770602354,74.03476173453267,666.312855610794,0.01527502736721497,0.01527502736721497),
    @("FAPs","Thbs1","Itga3","ECs",3,1,121.8208923333333,365.462677,0.4456653109566078,0.4456653109566078,3,1,12.77892533333333,38.336776,0.7206984557633653,0.7206984557633654,1556.740087167706,14010.66078450935,0.3211903013937273,0.3211903013937273),
    @("FAPs","Thbs1","Itga3","FAPs",3,1,121.8208923333333,365.462677,0.4456653109566078,0.4456653109566078,3,1,0.9898276666666668,2.969483,0.05582372947885773,0.05582372947885774,120.5816896095546,1085.235206485991,0.02487869975695269,0.02487869975695269),
    @("FAPs","Thbs1","Itga3","M2",3,1,121.8208923333333,365.462677,0.4456653109566078,0.4456653109566078,1,0.3333333333333333,0.02976266666666667,0.089288,0.001678537697541373,0.001678537697541373,3.62571461155289,32.631431503976,0.0007480660249271645,0.0007480660249271646),
    @("FAPs","Thbs1","Itga3","sCs",3,1,121.8208923333333,365.462677,0.4456653109566078,0.4456653109566078,3,1,3.932791,11.798373,0.2217992770602354,0.2217992770602354,479.0961089805024,4311.864980824522,0.09884824378100063,0.09884824378100063),
    @("M2","Thbs1","Itga3","ECs",3,1,87.673585,263.020755,0.3207419907481189,0.3207419907481188,3,1,12.77892533333333,38.336776,0.7206984557633653,0.7206984557633654,1120.374196420653,10083.36776778588,0.2311582574306369,0.2311582574306369),
    @("M2","Thbs1","Itga3","FAPs",3,1,87.673585,263.020755,0.3207419907481189,0.3207419907481188,3,1,0.9898276666666668,2.969483,0.05582372947885773,0.05582372947885774,86.78174006885168,781.0356606196651,0.01790501412403328,0.01790501412403328),
    @("M2","Thbs1","Itga3","M2",3,1,87.673585,263.020755,0.3207419907481189,0.3207419907481188,1,0.3333333333333333,0.02976266666666667,0.089288,0.001678537697541373,0.001678537697541373,2.609399685826667,23.48459717244,0.0005383775226551838,0.0005383775226551838),
    @("M2","Thbs1","Itga3","sCs",3,1,87.673585,263.020755,0.3207419907481189,0.3207419907481188,3,1,3.932791,11.798373,0.2217992770602354,0.2217992770602354,344.801886025735,3103.216974231615,0.07114034167079349,0.07114034167079349),
    @("sCs","Thbs1","Itga3","ECs",3,1,45.02666966666666,135.080009,0.1647240005714903,0.1647240005714903,3,1,12.77892533333333,38.336776,0.7206984557633653,0.7206984557633654,575.3924496789982,5178.532047110984,0.1187163328390368,0.1187163328390368),
    @("sCs","Thbs1","Itga3","FAPs",3,1,45.02666966666666,135.080009,0.1647240005714903,0.1647240005714903,3,1,0.9898276666666668,2.969483,0.05582372947885773,0.05582372947885774,44.56864337392744,401.117790365347,0.00919550804657808,0.00919550804657808),
    @("sCs","Thbs1","Itga3","M2",3,1,45.02666966666666,135.080009,0.1647240005714903,0.1647240005714903,1,0.3333333333333333,0.02976266666666667,0.089288,0.001678537697541373,0.001678537697541373,1.340113760399111,12.061023843592,0.0002764954446490731,0.0002764954446490731),
    @("sCs","Thbs1","Itga3","sCs",3,1,45.02666966666666,135.080009,0.1647240005714903,0.1647240005714903,3,1,3.932791,11.798373,0.2217992770602354,0.2217992770602354,177.0804812250396,1593.724331025357,0.03653566424122635,0.03653566424122635)
)

$startRow = 2
$numCols = 20
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $arr = New-Object 'object[,]' 1,$numCols
    for ($c = 0; $c -lt $numCols; $c++) {
        $arr[0,$c] = $data[$i][$c]
    }
    $rng = $ws.Range($ws.Cells.Item($r,1), $ws.Cells.Item($r,$numCols))
    $rng.Value = $arr
}
